$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A2:H17")
$key = $ws.Range("E2:E17")
$rng.Sort($key, 1)
